$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.468.73"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.802.54"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.58"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.74"
$ws.Range("E8").Value = "  +5.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.301"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0954"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.063.83"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.23"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.799.08"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.643"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.459.44"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.04"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0799"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.34"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "174.04"
$ws.Range("E24").Value = "  +3.70%  "
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.81"
$ws.Range("E26").Value = "  +6.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.82"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0531"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.395.15"
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.45"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("E41").Value = "  +2.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.949"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.62"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  +3.40%  "
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.963.50"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.99"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  +1.16%  "
